# Fruta / hortaliza, semanal
# Rows in the sheet (2,4,5,6,7,8) get their data re-ordered / updated to
# reflect the corrected weekly data. Row 3 is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row: D (Fecha serial), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# N (Unidad de comercializacion), P (Precio $/Kg), Q (Kg o Unidades)

$rows = @{
    2 = @{ D = 44377; I = "Segunda"; J = 550; K = 2000; L = 2800; M = 2364; N = "`$/docena de matas"; P = 394; Q = 6 }
    4 = @{ D = 45218; I = "Primera"; J = 180; K = 1400; L = 1500; M = 1444; N = "`$/docena de matas"; P = 241; Q = 6 }
    5 = @{ D = 45225; I = "Primera"; J = 60;  K = 1500; L = 2000; M = 1750; N = "`$/docena de matas"; P = 292; Q = 6 }
    6 = @{ D = 44267; I = "Primera"; J = 120; K = 1500; L = 1800; M = 1650; N = "`$/docena de matas"; P = 275; Q = 6 }
    7 = @{ D = 44623; I = "Primera"; J = 300; K = 1800; L = 2000; M = 1900; N = "`$/paquete";         P = 1900; Q = 1 }
    8 = @{ D = 44370; I = "Segunda"; J = 100; K = 1000; L = 1200; M = 1080; N = "`$/docena de matas"; P = 180; Q = 6 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("D$r").Value = $data.D
    $ws.Range("I$r").Value = $data.I
    $ws.Range("J$r").Value = $data.J
    $ws.Range("K$r").Value = $data.K
    $ws.Range("L$r").Value = $data.L
    $ws.Range("M$r").Value = $data.M
    $ws.Range("N$r").Value = $data.N
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
}
